$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 11500
$ws.Range("I40").Value = 11500
$ws.Range("K40").Value = 11500
$ws.Range("M40").Value = -11325

# Row 64
$ws.Range("H64").Value = 13499
$ws.Range("I64").Value = 13499
$ws.Range("K64").Value = 13499
$ws.Range("M64").Value = -13251

# Row 67
$ws.Range("H67").Value = 13499
$ws.Range("I67").Value = 13499
$ws.Range("K67").Value = 13499
$ws.Range("M67").Value = -12641

# Row 111
$ws.Range("H111").Value = 4773.963
$ws.Range("I111").Value = 4768.5
$ws.Range("K111").Value = 14305.5
$ws.Range("M111").Value = -11238.5

# Row 116
$ws.Range("H116").Value = 4418.2104
$ws.Range("I116").Value = 4327.364
$ws.Range("K116").Value = 4327.364
$ws.Range("M116").Value = -885.3639999999996

# Row 135
$ws.Range("H135").Value = 14707415
$ws.Range("I135").Value = 1136.9
$ws.Range("K135").Value = 10232.1
$ws.Range("M135").Value = -7697.1

# Row 138
$ws.Range("H138").Value = 10757977
$ws.Range("I138").Value = 1281.2778
$ws.Range("J138").Value = 25651864
$ws.Range("K138").Value = 3843.8334
$ws.Range("L138").Value = 76955592
$ws.Range("M138").Value = 1296.1666
$ws.Range("N138").Value = -76965872

# Row 141
$ws.Range("H141").Value = 2710.375
$ws.Range("I141").Value = 2710.375
$ws.Range("K141").Value = 8131.125
$ws.Range("M141").Value = -2951.125

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 4608
$ws.Range("I102").Value = 3676.3333
$ws.Range("K102").Value = 3676.3333
$ws.Range("M102").Value = -2054.3333

# Row 132
$ws.Range("H132").Value = 30305558
$ws.Range("I132").Value = 2290.5862
$ws.Range("K132").Value = 6871.758600000001
$ws.Range("M132").Value = -4341.758600000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2791.8125
$ws.Range("J20").Value = 1931
$ws.Range("L20").Value = 1931
$ws.Range("N20").Value = -2425

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 24395094
$ws.Range("I31").Value = 3948.074
$ws.Range("J31").Value = 71435160
$ws.Range("K31").Value = 3948.074
$ws.Range("L31").Value = 71435160
$ws.Range("M31").Value = -3653.074
$ws.Range("N31").Value = -71435750

# Row 34
$ws.Range("H34").Value = 24395094
$ws.Range("I34").Value = 3948.074
$ws.Range("J34").Value = 71435160
$ws.Range("K34").Value = 3948.074
$ws.Range("L34").Value = 71435160
$ws.Range("M34").Value = -3746.074
$ws.Range("N34").Value = -71435564

# Row 62
$ws.Range("H62").Value = 6619.273
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 6881.2
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 6881.2
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -8129.2

# Row 65
$ws.Range("H65").Value = 6619.273
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 6881.2
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 34406
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -40646

# Row 100
$ws.Range("H100").Value = 58333.332
$ws.Range("I100").Value = 50000
$ws.Range("J100").Value = 62500
$ws.Range("K100").Value = 50000
$ws.Range("L100").Value = 62500
$ws.Range("M100").Value = -48918
$ws.Range("N100").Value = -64664

# Row 132
$ws.Range("H132").Value = 4348.8945
$ws.Range("I132").Value = 3812.7222
$ws.Range("K132").Value = 11438.1666
$ws.Range("M132").Value = -8908.1666

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 249
$ws.Range("J22").Value = 50
$ws.Range("L22").Value = 150
$ws.Range("N22").Value = -488

# Row 27
$ws.Range("H27").Value = 249
$ws.Range("J27").Value = 50
$ws.Range("L27").Value = 150
$ws.Range("N27").Value = -354

# Row 109
$ws.Range("H109").Value = 3648.4443
$ws.Range("J109").Value = 2920
$ws.Range("L109").Value = 8760
$ws.Range("N109").Value = -10840

# Row 113
$ws.Range("H113").Value = 3725.75
$ws.Range("I113").Value = 2966.3333
$ws.Range("J113").Value = 3978.889
$ws.Range("K113").Value = 8898.999899999999
$ws.Range("L113").Value = 11936.667
$ws.Range("M113").Value = -6728.999899999999
$ws.Range("N113").Value = -16276.667

# Row 129
$ws.Range("H129").Value = 4107.8066
$ws.Range("I129").Value = 4061
$ws.Range("J129").Value = 4126.9546
$ws.Range("K129").Value = 12183
$ws.Range("L129").Value = 12380.8638
$ws.Range("M129").Value = -7183
$ws.Range("N129").Value = -22380.8638

# Row 131
$ws.Range("H131").Value = 32877.61
$ws.Range("J131").Value = 5214.5386
$ws.Range("L131").Value = 15643.6158
$ws.Range("N131").Value = -25723.6158

# Row 132
$ws.Range("H132").Value = 1669342.8
$ws.Range("I132").Value = 2259.75
$ws.Range("J132").Value = 3336425.8
$ws.Range("K132").Value = 20337.75
$ws.Range("L132").Value = 30027832.2
$ws.Range("M132").Value = -17807.75
$ws.Range("N132").Value = -30032892.2

# Row 134
$ws.Range("H134").Value = 2416.8
$ws.Range("I134").Value = 1491.3684
$ws.Range("K134").Value = 4474.1052
$ws.Range("M134").Value = 595.8948

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4091
$ws.Range("I70").Value = 3977
$ws.Range("K70").Value = 3977
$ws.Range("M70").Value = -3707

# Row 73
$ws.Range("H73").Value = 4091
$ws.Range("I73").Value = 3977
$ws.Range("K73").Value = 3977
$ws.Range("M73").Value = -3041

# Row 80
$ws.Range("H80").Value = 2740.7144
$ws.Range("I80").Value = 2780.8333
$ws.Range("K80").Value = 2780.8333
$ws.Range("M80").Value = -1782.8333

# Row 83
$ws.Range("H83").Value = 2740.7144
$ws.Range("I83").Value = 2780.8333
$ws.Range("K83").Value = 13904.1665
$ws.Range("M83").Value = -8912.166499999999

# Row 132
$ws.Range("H132").Value = 3531.5925
$ws.Range("I132").Value = 3645.8262
$ws.Range("K132").Value = 10937.4786
$ws.Range("M132").Value = -8407.4786

# Row 136
$ws.Range("H136").Value = 26173.936
$ws.Range("J136").Value = 26173.936
$ws.Range("L136").Value = 78521.808
$ws.Range("N136").Value = -83621.808

# Row 141
$ws.Range("H141").Value = 98000
$ws.Range("J141").Value = 105000
$ws.Range("L141").Value = 105000
$ws.Range("N141").Value = -115360

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4428.227
$ws.Range("I40").Value = 4004.8333
$ws.Range("J40").Value = 4936.3
$ws.Range("K40").Value = 4004.8333
$ws.Range("L40").Value = 4936.3
$ws.Range("M40").Value = -3868.8333
$ws.Range("N40").Value = -5208.3

# Row 135
$ws.Range("H135").Value = 40695
$ws.Range("J135").Value = 41000
$ws.Range("L135").Value = 41000
$ws.Range("N135").Value = -51140

# Row 141
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("WVR")
# Row 110
$ws.Range("H110").Value = 84330
$ws.Range("J110").Value = 84330
$ws.Range("L110").Value = 84330
$ws.Range("N110").Value = -92510

# Row 122
$ws.Range("H122").Value = 45455844
$ws.Range("I122").Value = 52632790
$ws.Range("J122").Value = 1835.6666
$ws.Range("K122").Value = 157898370
$ws.Range("L122").Value = 5506.9998
$ws.Range("M122").Value = -157895920
$ws.Range("N122").Value = -10406.9998

# Row 126
$ws.Range("H126").Value = 2901.6924
$ws.Range("I126").Value = 3386.0645
$ws.Range("K126").Value = 10158.1935
$ws.Range("M126").Value = -7688.193499999999

# Row 132
$ws.Range("H132").Value = 4354.163
$ws.Range("I132").Value = 4320.452
$ws.Range("J132").Value = 4556.4287
$ws.Range("K132").Value = 12961.356
$ws.Range("L132").Value = 13669.2861
$ws.Range("M132").Value = -10431.356
$ws.Range("N132").Value = -18729.2861

Write-Output "Updated 192 cells"